# COMMIT COMPLETE TEST FRAMEWORK
# Adds the "AddNewCustomer" worksheet (with sample customer rows and
# mailto: hyperlinks on the Email column) after "Sheet1", and updates
# Sheet1's selection.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$sheet1 = $wb.Worksheets.Item(1)

# --- Set Sheet1's selection to A1:B2 (it stays the inactive sheet once
# the new sheet below is added and becomes active) --------------------
[void]$sheet1.Range("A1:B2").Select()

# --- Force the "Hyperlink" cell style (font + cellStyleXf + cellXf) to
# be allocated *before* the "quotePrefix" style that DOB entries need
# below, by exercising Hyperlinks.Add on a disposable scratch sheet
# that is immediately removed again. -----------------------------------
$scratch = $wb.Worksheets.Add()
$scratch.Range("A1").Value = "scratch@example.com"
[void]$scratch.Hyperlinks.Add($scratch.Range("A1"), "mailto:scratch@example.com")
[void]$scratch.Delete()

# --- Add the new worksheet right after Sheet1 -------------------------
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "AddNewCustomer"

# --- Header row (C1:L1 first, then B1/A1) ------------------------------
$ws2.Range("C1").Value = "Name"
$ws2.Range("D1").Value = "Gender"
$ws2.Range("E1").Value = "DOB"
$ws2.Range("F1").Value = "Address"
$ws2.Range("G1").Value = "City"
$ws2.Range("H1").Value = "State"
$ws2.Range("I1").Value = "PinCode"
$ws2.Range("J1").Value = "TelephoneNo"
$ws2.Range("K1").Value = "Email"
$ws2.Range("L1").Value = "Password"
$ws2.Range("B1").Value = "LoginPWD"
$ws2.Range("A1").Value = "LoginUID"

# --- Row 2 values (first occurrences establish shared-string order) ---
$ws2.Range("D2").Value = "M"
$ws2.Range("F2").Value = "A4 PRIMECITY"
$ws2.Range("G2").Value = "BANGALORE"
$ws2.Range("H2").Value = "KARNATAKA"
$ws2.Range("L2").Value = "RAJ"

$ws2.Range("D3").Value = "F"

$ws2.Range("E2").Value = "'04071971"
$ws2.Range("E3").Value = "'04081972"

$ws2.Range("C2").Value = "RAJAN"
$ws2.Range("C3").Value = "BAJAN"

$ws2.Range("K2").Value = "RAJ111@GMAIL.COM"
$ws2.Range("K3").Value = "TAJ3@GMAIL.COM"

# --- Remaining cells (these reuse strings already interned above) -----
$ws2.Range("A2").Value = "mngr250914"
$ws2.Range("B2").Value = "rUzYdap"
$ws2.Range("A3").Value = "mngr250914"
$ws2.Range("B3").Value = "rUzYdap"
$ws2.Range("F3").Value = "A4 PRIMECITY"
$ws2.Range("G3").Value = "BANGALORE"
$ws2.Range("H3").Value = "KARNATAKA"
$ws2.Range("L3").Value = "RAJ"

$ws2.Range("I2").Value = 560100
$ws2.Range("J2").Value = 123456789
$ws2.Range("I3").Value = 560100
$ws2.Range("J3").Value = 213456789

# --- Hyperlinks on the Email column ------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("K2"), "mailto:RAJ111@GMAIL.COM")
$ws2.Hyperlinks.Add($ws2.Range("K3"), "mailto:TAJ3@GMAIL.COM")

# --- Column widths (best-fit approximations) ---------------------------
$ws2.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws2.Columns.Item(5).ColumnWidth = 8
$ws2.Columns.Item(6).ColumnWidth = 11.666666666666666
$ws2.Columns.Item(8).ColumnWidth = 10
$ws2.Columns.Item(10).ColumnWidth = 11.166666666666666
$ws2.Columns.Item(11).ColumnWidth = 18.333333333333332

# --- Selection / active cell on the new sheet ---------------------------
[void]$ws2.Range("E5").Select()
